$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, shifting existing rows 118:136 down to 119:137,
# then populate the new row with this week's price observation.
$ws.Rows.Item(118).Insert()

$ws.Range("A118").Value = 4
$ws.Range("B118").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C118").Value = "Los Lagos"
$ws.Range("D118").Value = 44474
$ws.Range("E118").Value = 10
$ws.Range("F118").Value = 100112017
$ws.Range("G118").Value = "Apio"
$ws.Range("H118").Value = "Americana (o)"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 50
$ws.Range("K118").Value = 11000
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = 11500
$ws.Range("N118").Value = "`$/docena de matas"
$ws.Range("O118").Value = "Región de Coquimbo"
$ws.Range("P118").Value = 1917
$ws.Range("Q118").Value = 6
$ws.Range("R118").Value = "Hortaliza"
